# B6-PowerPoint.pptx edit — Sat, Jun 13, 2020  9:05:42 PM
#
# The three data tables (on the slides that hold gridCol widths of
# 2879725/1547800/1547825, 3424250/1887525/1889125 and
# 2881325/1547800/1547825 EMU respectively — i.e. slides 14, 15 and 16)
# had their table style switched from the deck's custom "Table_0" style
# ({292A1B24-83BC-48D7-BAFD-8C12C31041D8}) to the built-in PowerPoint
# table style {25F10D69-1754-448F-978D-6932D95BDAFC}.

$p = $ppt.ActivePresentation

$newStyleId = "{25F10D69-1754-448F-978D-6932D95BDAFC}"
$targetSlides = 14, 15, 16

foreach ($slideIdx in $targetSlides) {
    $slide = $p.Slides.Item($slideIdx)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)

        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
